$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '74.925.24'
$ws.Range('E2').Value = '  +1.68%  '
$ws.Range('D3').Value = '2.817.36'
$ws.Range('E3').Value = '  +7.49%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '187.40'
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '595.03'
$ws.Range('E6').Value = '  +2.36%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +3.08%  '
$ws.Range('E9').Value = '  -4.53%  '
$ws.Range('D10').Value = '2.817.13'
$ws.Range('E10').Value = '  +7.52%  '
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('E12').Value = '  +3.64%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.88'
$ws.Range('E13').Value = '  +2.81%  '
$ws.Range('D14').Value = '3.341.39'
$ws.Range('E14').Value = '  +7.71%  '
$ws.Range('D15').Value = '74.962.65'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000188'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.81'
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('D18').Value = '2.821.64'
$ws.Range('E18').Value = '  +7.63%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '8.95'
$ws.Range('E19').Value = '  -1.10%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.30'
$ws.Range('E20').Value = '  +4.03%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '377.37'
$ws.Range('E21').Value = '  +1.81%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.26'
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.07%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '70.86'
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.961.93'
$ws.Range('E26').Value = '  +7.34%  '
$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '4.17'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.73'
$ws.Range('E28').Value = '  +4.04%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0000104'
$ws.Range('E29').Value = '  +10.61%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '516.22'
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.39'
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.71'
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.78'
$ws.Range('E34').Value = '  +2.69%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '162.39'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '19.93'
$ws.Range('E37').Value = '  +4.08%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.118'
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '19.37'
$ws.Range('E39').Value = '  +0.66%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '186.85'
$ws.Range('E40').Value = '  +16.30%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.340'
$ws.Range('E42').Value = '  +3.42%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.00'
$ws.Range('E43').Value = '  +2.09%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.67'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.21'
$ws.Range('E45').Value = '  +2.18%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '40.00'
$ws.Range('E46').Value = '  +2.72%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.33'
$ws.Range('E47').Value = '  -1.53%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0852'
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.574'
$ws.Range('E49').Value = '  +8.77%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.72'
$ws.Range('E50').Value = '  +2.63%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.634'
$ws.Range('E51').Value = '  +8.41%  '
